# The paragraph in question renders (as visible text) like:
#   On les souloict donner avecq <tl>manches</tl> ou <tl>poches de <m>cuir</m></tl>, qui
# where "<tl>", "</tl>", "<m>" etc. are literal characters living in their own
# (blue, Courier New) runs, interleaved with plain-black runs holding the
# actual words. The edit collapses:
#   <tl>manches</tl> ou <tl>poches de
# down to just:
#   <tl>manches ou poches de
# i.e. the run holding "manches" keeps its own (black) formatting but grows
# to read "manches ou poches de ", while the runs for "</tl>", " ou ",
# "<tl>" and "poches de " are removed outright. The following "<m>cuir</m>"
# (etc.) runs are left completely untouched.

$d = $word.ActiveDocument

# Locate the "manches" run via Find - Execute() collapses/moves the range
# onto the matched text, so $wordRange ends up exactly spanning that run
# (and nothing else), preserving its own (black) run formatting when we
# rewrite its text below.
$wordRange = $d.Content
$found = $wordRange.Find.Execute("manches", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'manches' run to edit"
}

# Grow this run's text in place (keeps its black-colored formatting).
$wordRange.Text = "manches ou poches de "

# Find where the following "<m>" marker run begins, so we know exactly how
# much trailing text (the now-redundant "</tl> ou <tl>poches de " runs) to
# strip out before it.
$tailRange = $d.Range($wordRange.End, $d.Content.End)
$foundTail = $tailRange.Find.Execute("<m>", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)

if (-not $foundTail) {
    throw "Could not find the trailing '<m>' marker"
}

$deadRange = $d.Range($wordRange.End, $tailRange.Start)
$deadRange.Delete()
